$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing requisito rows 44-46 (text changes only; keep position/formatting as-is)
$ws.Range("B44").Value = "LOQ4010 -  Introdução à  Engenharia  Química  (Requisito)`n"
$ws.Range("C44").Value = "LOQ4010 -  Introdução à  Engenharia  Química  (Requisito)`n"

$ws.Range("B45").Value = "LOQ4095 -  Química Geral Experimental  (Requisito)`n"
$ws.Range("C45").Value = "LOQ4095 -  Química Geral Experimental  (Requisito)`n"

$ws.Range("B46").Value = "LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito)`n"
$ws.Range("C46").Value = "LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito)`n"

# Add new row 47 with the same style/formatting as row 46, then set its text
$ws.Range("B46:C46").Copy()
$ws.Range("B47:C47").PasteSpecial(-4122)
$ws.Range("B47").Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)`n"
$ws.Range("C47").Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)`n"
$ws.Rows.Item(47).RowHeight = 30
